$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "ProgramsTab" SQL query text in B2 ---
# Adds a CASE expression for the "Website" column (program_link / program_acronym)
# and tweaks whitespace/formatting around the existing "Data Location Details"
# CASE block and the FROM/WHERE clause, matching the author's edit.
$newQuery = @"
SELECT DISTINCT 
    prg.program_name AS "Program",
  CASE
    WHEN prg.program_link IS NOT NULL THEN prg.program_acronym
        ELSE prg.program_link
    END  AS "Website",
    prg.focus_area AS "Focus Area",
    prg.cancer_type AS "Cancer Type",
 CASE 
        WHEN prg.data_link IS NOT NULL THEN prg.website       
        ELSE prg.data_link
    END AS "Data Location Details"
FROM 
    df_program prgWHERE 
     prg.cancer_type LIKE '%Kidney Cancer%'
ORDER BY 
    lower(prg.program_name) ASC
LIMIT 100;
"@

$ws.Range("B2").Value2 = $newQuery

# --- Update the sheet view: scroll so row 5 is at the top, and move the
#     active selection down to B20 (matching the saved view state) ---
$excel.ActiveWindow.ScrollRow = 5
$excel.ActiveWindow.ScrollColumn = 1
[void]$ws.Range("B20").Select()
